$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.312.89"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").Value = "2.423.87"

$ws.Range("D5").Value = "'556.14"
$ws.Range("E5").Value = "  +2.22%  "

$ws.Range("D6").Value = "'143.68"
$ws.Range("E6").Value = "  +4.97%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +1.72%  "

$ws.Range("D9").Value = "2.423.58"
$ws.Range("E9").Value = "  +3.21%  "

$ws.Range("E10").Value = "  +4.61%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("D14").Value = "'26.27"
$ws.Range("E14").Value = "  +6.33%  "

$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "  +9.27%  "

$ws.Range("D16").Value = "2.862.00"
$ws.Range("E16").Value = "  +3.17%  "

$ws.Range("D17").Value = "62.244.32"
$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("D18").Value = "2.424.40"
$ws.Range("E18").Value = "  +3.29%  "

$ws.Range("D19").Value = "'11.08"
$ws.Range("E19").Value = "  +4.16%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  +1.68%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'324.50"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").Value = "'6.73"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  +5.83%  "

$ws.Range("D25").Value = "'64.80"
$ws.Range("E25").Value = "  +2.32%  "

$ws.Range("D26").Value = "'9.11"
$ws.Range("E26").Value = "  +9.30%  "

$ws.Range("D27").Value = "'572.90"
$ws.Range("E27").Value = "  +15.07%  "

$ws.Range("D28").Value = "2.538.63"
$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "'8.39"
$ws.Range("E30").Value = "  +5.42%  "

$ws.Range("D31").Value = "0.0₃0939"
$ws.Range("E31").Value = "  +9.35%  "

$ws.Range("E32").Value = "  +6.03%  "

$ws.Range("E33").Value = "  +2.05%  "

$ws.Range("E34").Value = "  +3.84%  "

$ws.Range("E35").Value = "  +5.13%  "

$ws.Range("D36").Value = "'5.72"
$ws.Range("E36").Value = "  +9.12%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.83"
$ws.Range("E38").Value = "  +5.09%  "

$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'18.77"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("D42").Value = "'149.61"
$ws.Range("E42").Value = "  +4.26%  "

$ws.Range("D44").Value = "'41.68"
$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = "  +14.07%  "

$ws.Range("D46").Value = "'150.95"
$ws.Range("E46").Value = "  +5.27%  "

$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("D48").Value = "'0.0544"
$ws.Range("E48").Value = "  +4.98%  "

$ws.Range("D49").Value = "'20.36"
$ws.Range("E49").Value = "  +6.80%  "

$ws.Range("E50").Value = "  +3.81%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  +3.54%  "
